# Updated cryptos list on Mon May  1 04:28:33 UTC 2023 with GitHub Actions
#
# Refreshes the Price (D) and Volume(1h) (E) columns for every coin row,
# and re-ranks two pairs of coins whose relative ordering changed:
#   - RenderToken now ranks above WEMIXToken (rows 44/45 swap)
#   - EnergySwap now ranks above Decentraland (rows 46/47 swap)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.583.01"
$ws.Range("E2").Value = "  -3.16%  "
$ws.Range("D3").Value = "1.849.29"
$ws.Range("E3").Value = "  -3.58%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  -0.98%  "
$ws.Range("D5").Value = "'335.68"
$ws.Range("E5").Value = "  +3.17%  "
$ws.Range("E6").Value = "  -0.98%  "
$ws.Range("D7").Value = "'0.4662"
$ws.Range("E7").Value = "  -3.10%  "
$ws.Range("D8").Value = "'0.3912"
$ws.Range("D9").Value = "'46.61"
$ws.Range("E9").Value = "  -2.79%  "
$ws.Range("D10").Value = "'0.07882"
$ws.Range("E10").Value = "  -4.04%  "
$ws.Range("D11").Value = "'0.9818"
$ws.Range("E11").Value = "  -2.69%  "
$ws.Range("D12").Value = "'22.17"
$ws.Range("E12").Value = "  -5.30%  "
$ws.Range("D13").Value = "1.879.22"
$ws.Range("E13").Value = "  -1.52%  "
$ws.Range("D14").Value = "'5.842"
$ws.Range("E14").Value = "  -3.43%  "
$ws.Range("D15").Value = "'7.005"
$ws.Range("E15").Value = "  -3.13%  "
$ws.Range("D16").Value = "'0.06834"
$ws.Range("E16").Value = "  -0.69%  "
$ws.Range("E17").Value = "  -1.02%  "
$ws.Range("D18").Value = "'87.60"
$ws.Range("D19").Value = "'0.00001011"
$ws.Range("E19").Value = "  -2.66%  "
$ws.Range("E20").Value = "  -2.94%  "
$ws.Range("D21").Value = "'1.001"
$ws.Range("E21").Value = "  -0.97%  "
$ws.Range("D22").Value = "28.596.96"
$ws.Range("E22").Value = "  -3.12%  "
$ws.Range("D23").Value = "'5.392"
$ws.Range("E23").Value = "  -4.88%  "
$ws.Range("D24").Value = "'11.24"
$ws.Range("E24").Value = "  -5.24%  "
$ws.Range("E25").Value = "  -3.16%  "
$ws.Range("D26").Value = "2.137.52"
$ws.Range("E26").Value = "  +0.21%  "
$ws.Range("D27").Value = "'153.21"
$ws.Range("E27").Value = "  -1.74%  "
$ws.Range("D28").Value = "'6.191"
$ws.Range("E28").Value = "  -5.26%  "
$ws.Range("E29").Value = "  -3.16%  "
$ws.Range("D30").Value = "'2.019"
$ws.Range("E30").Value = "  -3.67%  "
$ws.Range("D31").Value = "'117.39"
$ws.Range("E31").Value = "  -2.67%  "
$ws.Range("D32").Value = "'0.9745"
$ws.Range("E32").Value = "  -4.05%  "
$ws.Range("D33").Value = "'0.09448"
$ws.Range("E33").Value = "  -1.94%  "
$ws.Range("D34").Value = "'5.369"
$ws.Range("E34").Value = "  -4.41%  "
$ws.Range("E35").Value = "  -1.58%  "
$ws.Range("D36").Value = "'1.348"
$ws.Range("E36").Value = "  -1.73%  "
$ws.Range("D37").Value = "'0.06143"
$ws.Range("E37").Value = "  -2.57%  "
$ws.Range("D38").Value = "'0.02190"
$ws.Range("E38").Value = "  -4.10%  "
$ws.Range("D39").Value = "'1.162"
$ws.Range("E39").Value = "  -1.84%  "
$ws.Range("D40").Value = "'0.5686"
$ws.Range("E40").Value = "  -4.15%  "
$ws.Range("D41").Value = "'7.563"
$ws.Range("E41").Value = "  -4.38%  "
$ws.Range("D42").Value = "'10.10"
$ws.Range("E42").Value = "  -5.47%  "
$ws.Range("D43").Value = "'0.1789"
$ws.Range("E43").Value = "  -3.13%  "
# Rows 44-47 re-rank: RenderToken moves above WEMIXToken, and EnergySwap
# moves above Decentraland. Write the full new row contents directly
# (Coin name, Link, Price, Volume) rather than relying on a cell-swap, since
# the Price/Volume values were also refreshed.
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Value = "'2.379"
$ws.Range("E44").Value = "  -3.44%  "

$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").Value = "'1.250"
$ws.Range("E45").Value = "  -2.36%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'11.89"
$ws.Range("E46").Value = "  -4.22%  "

$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "'0.5376"
$ws.Range("E47").Value = "  -3.25%  "

$ws.Range("D48").Value = "'0.07141"
$ws.Range("E48").Value = "  -4.42%  "
$ws.Range("D49").Value = "'1.904"
$ws.Range("E49").Value = "  -1.82%  "
$ws.Range("D50").Value = "'113.12"
$ws.Range("E50").Value = "  -4.27%  "
$ws.Range("D51").Value = "'43.54"
$ws.Range("E51").Value = "  +3.33%  "
